# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
# Remove all event rows (2-14), leaving only the header row.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("2:14").Delete()

# --- Sheet 2: 演出 ---
# "想去人数" (want-to-go count) resets to 0 for every event.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,6).Value = 0
$ws2.Cells.Item(3,6).Value = 0
$ws2.Cells.Item(4,6).Value = 0
$ws2.Cells.Item(5,6).Value = 0

# --- Sheet 3: 本地生活 ---
# No changes (already header-only).

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)

# "想去人数" resets to 0 for every existing event row.
for ($r = 2; $r -le 17; $r++) {
    $ws4.Cells.Item($r, 6).Value = 0
}

# A new event row (duplicate of the "第二届北极光动漫展" entry, row 17) is
# inserted after row 17, pushing the "万圣漫控嘉年华10" row from 18 down to 19.
$ws4.Rows("17:17").Copy()
$ws4.Rows("18:18").Insert()

# Fix the sequence number / counter column for the new and shifted rows.
$ws4.Cells.Item(18,1).Value = 17
$ws4.Cells.Item(18,6).Value = 0
$ws4.Cells.Item(19,1).Value = 18
$ws4.Cells.Item(19,6).Value = 0

# Restore the plain (non-bordered) number style on the new row's A cell to
# match the rest of the column.
$ws4.Cells.Item(19,1).Copy()
$ws4.Cells.Item(18,1).PasteSpecial(-4122)
